# Automatische test-sync: 2025-06-24 20:14:50
# Adds a new "Verzoek om factuur" log entry to the Logs sheet and refreshes
# the Dashboard category counts to reflect it.

$wb = $excel.ActiveWorkbook

# ---- Logs sheet: append the new row ----
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A20").Value = "Verzoek om factuur"
$logs.Range("B20").Value = "mailmind.test@zohomail.eu"
$logs.Range("C20").Value = "Kunt u mij een factuur sturen voor mijn laatste bestelling?"
$logs.Range("D20").Value = "Factuur / Administratie"
$logs.Range("E20").Value = "Beste klant,`nDank u voor uw e-mail. Om u zo snel mogelijk te kunnen helpen, hebben we wat extra informatie nodig. Kunt u alstublieft uw klantgegevens (zoals uw naam en e-mailadres) en/of het factuurnummer van uw laatste bestelling doorgeven? Op die manier kunnen we de factuur voor u zo spoedig mogelijk opstellen en toesturen.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F20").Value = "2025-06-24 20:14:48"
$logs.Range("G20").Value = "Ja"

# ---- Logs sheet: extend the conditional formatting to cover the new row ----
$catFcs = $logs.Range("D2:D19").FormatConditions
$newCatRange = $logs.Range("D2:D20")
for ($i = 1; $i -le $catFcs.Count; $i++) {
    $catFcs.Item($i).ModifyAppliesToRange($newCatRange)
}

$answeredFcs = $logs.Range("G2:G19").FormatConditions
$newAnsweredRange = $logs.Range("G2:G20")
for ($i = 1; $i -le $answeredFcs.Count; $i++) {
    $answeredFcs.Item($i).ModifyAppliesToRange($newAnsweredRange)
}

# ---- Dashboard sheet: refresh the category counts ----
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A3").Value = "Factuur / Administratie"
$dash.Range("B3").Value = 3

$dash.Range("A4").Value = "IT / Technisch probleem"
$dash.Range("B4").Value = 3

$dash.Range("A5").Value = "Sollicitatie / Vacature"
$dash.Range("B5").Value = 2
